# Add a new "path" row right after the header row (becomes the new row 2),
# shifting all the existing data rows down by one, then stamp it with its
# status code and a distinguishing monospace font - mirrors the commit's
# "modified status code test data" change to the CTS data provider sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 2 - everything below (old rows 2-25) shifts to 3-26.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the new path/statusCode pair.
$ws.Range("A2").Value = "/Common/PopUps/popDefinition.aspx?id=485395&version=healthprofessional&language=English&dictionary=genetic"
$ws.Range("B2").Value = 200

# Give the new path cell its own look: Menlo 11pt in dark gray (#222222).
$ws.Range("A2").Font.Name = "Menlo"
$ws.Range("A2").Font.Size = 11
$ws.Range("A2").Font.Color = 2236962

# Leave the selection on the status-code cell of the new row.
$ws.Range("B2").Select()
